$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number + report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  10"
$ws.Range("C9").Value = "Report Covering the Week  3/3/2025  Through  3/9/2025"

# --- Reference cell (already formatted as style 13 text placeholder) used to restore formatting ---
# A14 holds a text "s13" placeholder cell; we copy-paste-special its format (no value change)
# onto cells that must become text placeholders, after writing their text value.

# --- Cells changing from numeric to text placeholder ("0" / "***.*") ---
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("C16").PasteSpecial(-4122)

$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("C22").PasteSpecial(-4122)

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("D27").PasteSpecial(-4122)

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("E27").PasteSpecial(-4122)

$ws.Range("F29").NumberFormat = "@"
$ws.Range("F29").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("F29").PasteSpecial(-4122)

$ws.Range("F30").NumberFormat = "@"
$ws.Range("F30").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("F30").PasteSpecial(-4122)

$ws.Range("F31").NumberFormat = "@"
$ws.Range("F31").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("F31").PasteSpecial(-4122)

# --- Cells changing from text placeholder to numeric (style must become 14, #,##0) ---
# I16 is a stable style-14 (#,##0) numeric cell that is never itself retyped, used as the
# format donor here.
$ws.Range("I16").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = 1

$ws.Range("I16").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 1

# --- Plain numeric value updates ---
$ws.Range("I15").Value = 3
$ws.Range("K15").Value = 0
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = -11.111111111111
$ws.Range("J16").Value = 20
$ws.Range("K16").Value = -20
$ws.Range("L16").Value = -33.333333333333
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = -86.324786324786
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = -22.222222222222
$ws.Range("I17").Value = 23
$ws.Range("J17").Value = 16
$ws.Range("K17").Value = 43.75
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -14.814814814814
$ws.Range("N17").Value = -54
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = 50
$ws.Range("I18").Value = 19
$ws.Range("J18").Value = 19
$ws.Range("L18").Value = 5.555555555555
$ws.Range("M18").Value = -26.923076923076
$ws.Range("N18").Value = -76.829268292682
$ws.Range("C19").Value = 17
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 88.888888888888
$ws.Range("F19").Value = 50
$ws.Range("G19").Value = 56
$ws.Range("H19").Value = -10.714285714285
$ws.Range("I19").Value = 112
$ws.Range("J19").Value = 119
$ws.Range("K19").Value = -5.882352941176
$ws.Range("L19").Value = -15.151515151515
$ws.Range("M19").Value = -8.196721311475
$ws.Range("N19").Value = -28.205128205128
$ws.Range("G20").Value = 2
$ws.Range("M20").Value = -66.666666666666
$ws.Range("N20").Value = -98
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 15
$ws.Range("E21").Value = 40
$ws.Range("F21").Value = 78
$ws.Range("G21").Value = 84
$ws.Range("H21").Value = -7.142857142857
$ws.Range("I21").Value = 175
$ws.Range("J21").Value = 188
$ws.Range("K21").Value = -6.914893617021
$ws.Range("L21").Value = -13.793103448275
$ws.Range("M21").Value = -11.167512690355
$ws.Range("N21").Value = -65.753424657534
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = -100
$ws.Range("G22").Value = 7
$ws.Range("H22").Value = -71.428571428571
$ws.Range("I22").Value = 5
$ws.Range("J22").Value = 9
$ws.Range("K22").Value = -44.444444444444
$ws.Range("L22").Value = -37.5
$ws.Range("M22").Value = 66.666666666666
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 11
$ws.Range("J23").Value = 8
$ws.Range("K23").Value = 37.5
$ws.Range("L23").Value = 57.142857142857
$ws.Range("M23").Value = 57.142857142857
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 10
$ws.Range("E24").Value = 100
$ws.Range("F24").Value = 48
$ws.Range("H24").Value = 14.285714285714
$ws.Range("I24").Value = 138
$ws.Range("J24").Value = 109
$ws.Range("K24").Value = 26.605504587156
$ws.Range("L24").Value = 7.8125
$ws.Range("M24").Value = -0.719424460431
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 160
$ws.Range("F25").Value = 37
$ws.Range("G25").Value = 21
$ws.Range("H25").Value = 76.190476190476
$ws.Range("I25").Value = 74
$ws.Range("J25").Value = 50
$ws.Range("K25").Value = 48
$ws.Range("L25").Value = 0
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 25
$ws.Range("F26").Value = 21
$ws.Range("G26").Value = 23
$ws.Range("H26").Value = -8.695652173913
$ws.Range("I26").Value = 47
$ws.Range("J26").Value = 65
$ws.Range("K26").Value = -27.692307692307
$ws.Range("L26").Value = -33.802816901408
$ws.Range("M26").Value = -14.545454545454
$ws.Range("I27").Value = 3
$ws.Range("K27").Value = -40
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 8
$ws.Range("J28").Value = 10
$ws.Range("K28").Value = -20
$ws.Range("L28").Value = -27.272727272727
